# Auto-generated edit script applying the cryptos.xlsx data refresh diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{Row=2; B=$null; C=$null; D="20.574.89"; E="  +2.37%  "},
    @{Row=3; B=$null; C=$null; D="1.473.31"; E="  +3.13%  "},
    @{Row=4; B=$null; C=$null; D="1.011"; E="  +0.98%  "},
    @{Row=5; B=$null; C=$null; D="0.9881"; E="  -1.37%  "},
    @{Row=6; B=$null; C=$null; D="281.52"; E="  +2.42%  "},
    @{Row=7; B=$null; C=$null; D="0.3726"; E="  +0.09%  "},
    @{Row=8; B=$null; C=$null; D="0.3204"; E="  +4.01%  "},
    @{Row=9; B=$null; C=$null; D="41.73"; E="  +3.82%  "},
    @{Row=10; B=$null; C=$null; D="1.071"; E="  +6.58%  "},
    @{Row=11; B=$null; C=$null; D="0.06728"; E="  +2.13%  "},
    @{Row=12; B=$null; C=$null; D="1.005"; E="  +0.32%  "},
    @{Row=13; B=$null; C=$null; D="5.654"; E="  +4.32%  "},
    @{Row=14; B=$null; C=$null; D="18.54"; E="  +7.44%  "},
    @{Row=15; B=$null; C=$null; D="6.323"; E="  +2.28%  "},
    @{Row=16; B=$null; C=$null; D="1.476.95"; E="  +3.10%  "},
    @{Row=17; B=$null; C=$null; D="0.00001041"; E="  +2.87%  "},
    @{Row=18; B=$null; C=$null; D="0.05809"; E="  -0.26%  "},
    @{Row=19; B=$null; C=$null; D="73.15"; E="  -3.19%  "},
    @{Row=20; B=$null; C=$null; D="0.9871"; E="  -1.44%  "},
    @{Row=21; B=$null; C=$null; D="5.728"; E="  +0.66%  "},
    @{Row=22; B=$null; C=$null; D="14.93"; E="  +2.85%  "},
    @{Row=23; B=$null; C=$null; D="11.26"; E="  +1.67%  "},
    @{Row=24; B=$null; C=$null; D="2.304"; E="  -1.29%  "},
    @{Row=25; B=$null; C=$null; D="20.699.26"; E="  +2.95%  "},
    @{Row=26; B=$null; C=$null; D="2.346"; E="  +2.46%  "},
    @{Row=27; B=$null; C=$null; D="138.25"; E="  -0.37%  "},
    @{Row=28; B=$null; C=$null; D="17.70"; E="  +4.69%  "},
    @{Row=29; B=$null; C=$null; D="1.645.32"; E="  +3.18%  "},
    @{Row=30; B=$null; C=$null; D="113.81"; E="  +3.90%  "},
    @{Row=31; B=$null; C=$null; D="3.976"; E="  +1.39%  "},
    @{Row=32; B=$null; C=$null; D="5.397"; E="  -0.69%  "},
    @{Row=33; B=$null; C=$null; D="0.8487"; E="  -7.18%  "},
    @{Row=34; B=$null; C=$null; D="1.642"; E="  +26.65%  "},
    @{Row=35; B=$null; C=$null; D="0.07873"; E="  +0.99%  "},
    @{Row=36; B=$null; C=$null; D="0.06112"; E="  +7.38%  "},
    @{Row=37; B=$null; C=$null; D="4.951"; E="  +3.60%  "},
    @{Row=38; B=$null; C=$null; D="10.81"; E="  -5.96%  "},
    @{Row=39; B="Frax"; C="https://coinranking.com/coin/KfWtaeV1W+frax-frax"; D="0.9920"; E="  -0.91%  "},
    @{Row=40; B="VeChain"; C="https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"; D="0.02085"; E="  +3.30%  "},
    @{Row=41; B="FraxShare"; C="https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"; D="7.784"; E="  -7.48%  "},
    @{Row=42; B=$null; C=$null; D="1.140"; E="  +1.46%  "},
    @{Row=43; B=$null; C=$null; D="0.1913"; E="  -0.37%  "},
    @{Row=44; B=$null; C=$null; D="0.5465"; E="  +2.38%  "},
    @{Row=45; B=$null; C=$null; D="12.65"; E="  +3.10%  "},
    @{Row=46; B=$null; C=$null; D=$null; E="  +1.20%  "},
    @{Row=47; B=$null; C=$null; D="121.61"; E="  +9.84%  "},
    @{Row=48; B=$null; C=$null; D="0.5390"; E="  +4.73%  "},
    @{Row=49; B=$null; C=$null; D="1.846"; E="  +3.89%  "},
    @{Row=50; B="Cronos"; C="https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"; D="0.06457"; E="  +4.22%  "},
    @{Row=51; B="EOS"; C="https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"; D="1.060"; E="  +0.50%  "},
)

foreach ($item in $updates) {
    if ($item.B -ne $null) {
        $ws.Cells.Item($item.Row, 2).Value = $item.B
    }
    if ($item.C -ne $null) {
        $ws.Cells.Item($item.Row, 3).Value = $item.C
    }
    if ($item.D -ne $null) {
        $dcell = $ws.Cells.Item($item.Row, 4)
        $dcell.NumberFormat = "@"
        $dcell.Value = $item.D
        $dcell.Style = "Normal"
    }
    if ($item.E -ne $null) {
        $ws.Cells.Item($item.Row, 5).Value = $item.E
    }
}

